$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.894.88"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "3.782.59"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'603.47"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").Value = "'163.50"
$ws.Range("E6").Value = "  -1.70%  "
$ws.Range("D7").Value = "3.780.16"
$ws.Range("E7").Value = "  -0.96%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").Value = "'6.80"
$ws.Range("E12").Value = "  +8.01%  "
$ws.Range("D13").Value = "'0.0000246"
$ws.Range("E13").Value = "  -2.73%  "
$ws.Range("D14").Value = "'35.06"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").Value = "4.416.11"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("D16").Value = "3.774.92"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "67.863.44"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "'18.18"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("E19").Value = "  +1.90%  "
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("D21").Value = "'458.03"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").Value = "'9.46"
$ws.Range("E22").Value = "  -4.25%  "
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("D24").Value = "'83.16"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("E25").Value = "  -2.53%  "
$ws.Range("E26").Value = "  -1.74%  "
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "'9.89"
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("D30").Value = "3.933.84"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("E31").Value = "  -6.77%  "
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("E33").Value = "  -1.89%  "
$ws.Range("E34").Value = "  -2.02%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E36").Value = "  -1.83%  "
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("D38").Value = "'0.148"
$ws.Range("E38").Value = "  +7.31%  "
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").Value = "'3.20"
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("D41").Value = "'0.976"
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D44").Value = "'43.64"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").Value = "'47.05"
$ws.Range("E45").Value = "  -2.17%  "
$ws.Range("D46").Value = "'152.58"
$ws.Range("E46").Value = "  +2.58%  "
$ws.Range("D47").Value = "'0.293"
$ws.Range("E47").Value = "  -2.20%  "
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("D51").Value = "'26.58"
$ws.Range("E51").Value = "  -6.91%  "
